# Recompute of the NATMI LR-pair output (Penk-Oprd1) with refreshed TPM values.
# The sending/target cluster set changed (Inflammatory-Mac & MuSCs dropped,
# Resolving-Mac added as a target cluster), and all expression/specificity
# statistics were recalculated accordingly, shrinking the table from 5 to 4 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs (Penk/Oprd1)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Penk"
$ws.Range("C2").Value = "Oprd1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1581486666666667
$ws.Range("H2").Value = 0.474446
$ws.Range("I2").Value = 0.005396672170837162
$ws.Range("J2").Value = 0.005396672170837161
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02165966666666666
$ws.Range("N2").Value = 0.064979
$ws.Range("O2").Value = 0.7290280598220596
$ws.Range("P2").Value = 0.7290280598220598
$ws.Range("Q2").Value = 0.003425447403777777
$ws.Range("R2").Value = 0.030829026634
$ws.Range("S2").Value = 0.003934325442201119
$ws.Range("T2").Value = 0.003934325442201119

# Row 3: ECs -> Resolving-Mac (Penk/Oprd1)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Penk"
$ws.Range("C3").Value = "Oprd1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1581486666666667
$ws.Range("H3").Value = 0.474446
$ws.Range("I3").Value = 0.005396672170837162
$ws.Range("J3").Value = 0.005396672170837161
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.008050666666666666
$ws.Range("N3").Value = 0.024152
$ws.Range("O3").Value = 0.2709719401779404
$ws.Range("P3").Value = 0.2709719401779404
$ws.Range("Q3").Value = 0.001273202199111111
$ws.Range("R3").Value = 0.011458819792
$ws.Range("S3").Value = 0.001462346728636043
$ws.Range("T3").Value = 0.001462346728636043

# Row 4: FAPs -> FAPs (Penk/Oprd1)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Penk"
$ws.Range("C4").Value = "Oprd1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 29.146701
$ws.Range("H4").Value = 87.44010299999999
$ws.Range("I4").Value = 0.9946033278291628
$ws.Range("J4").Value = 0.9946033278291628
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02165966666666666
$ws.Range("N4").Value = 0.064979
$ws.Range("O4").Value = 0.7290280598220596
$ws.Range("P4").Value = 0.7290280598220598
$ws.Range("Q4").Value = 0.6313078280929999
$ws.Range("R4").Value = 5.681770452836999
$ws.Range("S4").Value = 0.7250937343798586
$ws.Range("T4").Value = 0.7250937343798587

# Row 5: FAPs -> Resolving-Mac (Penk/Oprd1)
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Penk"
$ws.Range("C5").Value = "Oprd1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.146701
$ws.Range("H5").Value = 87.44010299999999
$ws.Range("I5").Value = 0.9946033278291628
$ws.Range("J5").Value = 0.9946033278291628
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008050666666666666
$ws.Range("N5").Value = 0.024152
$ws.Range("O5").Value = 0.2709719401779404
$ws.Range("P5").Value = 0.2709719401779404
$ws.Range("Q5").Value = 0.2346503741839999
$ws.Range("R5").Value = 2.111853367656
$ws.Range("S5").Value = 0.2695095934493043
$ws.Range("T5").Value = 0.2695095934493043

# The old row 6 (Resolving-Mac -> FAPs) combination no longer exists under the new
# cluster set, so the trailing row is removed and the sheet dimension shrinks to T5.
$ws.Rows.Item(6).Delete()
